$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Wrap the second ColorHunt palette URL (paragraph 8) in a hyperlink,
#    same as the first palette URL a paragraph above it.
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
$r8.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark from the range
$d.Hyperlinks.Add($r8, "https://colorhunt.co/palette/f2debaffefd60e5e6f3a8891") | Out-Null

# ---------------------------------------------------------------------------
# Helper: appends a new paragraph at the very end of the document containing
# the given plain text (can be empty for a blank paragraph) and returns the
# Paragraph object that was just created.
# ---------------------------------------------------------------------------
function Add-EndParagraph([string]$text) {
    $r = $d.Content
    $r.Collapse(0) | Out-Null
    $r.InsertParagraphAfter() | Out-Null
    $r.Collapse(0) | Out-Null
    if ($text -ne "") {
        $r.InsertAfter($text) | Out-Null
    }
    return $d.Paragraphs.Last
}

# Helper: within a paragraph's range, find the given text and turn it into a
# hyperlink pointing at the given address.
function Add-InlineHyperlink($paragraph, [string]$textToFind, [string]$address) {
    $f = $paragraph.Range.Duplicate
    $f.Find.Execute($textToFind) | Out-Null
    $d.Hyperlinks.Add($f, $address) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Blank paragraph separating the color palettes from the new image
#    credits.
# ---------------------------------------------------------------------------
Add-EndParagraph "" | Out-Null

# ---------------------------------------------------------------------------
# 3. img7 - Shahadat Rahman on Unsplash
# ---------------------------------------------------------------------------
$pImg7 = Add-EndParagraph "img7: Photo by Shahadat Rahman on Unsplash"
Add-InlineHyperlink $pImg7 "Shahadat Rahman" "https://unsplash.com/@shahadatrhidoy?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"
Add-InlineHyperlink $pImg7 "Unsplash" "https://unsplash.com/photos/HL8dH6PxOQI?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"

# ---------------------------------------------------------------------------
# 4. img8 - fabio on Unsplash
# ---------------------------------------------------------------------------
$pImg8 = Add-EndParagraph "img8: Photo by fabio on Unsplash"
Add-InlineHyperlink $pImg8 "fabio" "https://unsplash.com/@fabioha?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"
Add-InlineHyperlink $pImg8 "Unsplash" "https://unsplash.com/photos/oyXis2kALVg?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"

# ---------------------------------------------------------------------------
# 5. img9 - Kevin Ku on Unsplash
# ---------------------------------------------------------------------------
$pImg9 = Add-EndParagraph "img9: Photo by Kevin Ku on Unsplash"
Add-InlineHyperlink $pImg9 "Kevin Ku" "https://unsplash.com/@ikukevk?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"
Add-InlineHyperlink $pImg9 "Unsplash" "https://unsplash.com/photos/w7ZyuGYNpRQ?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
